$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.059.19"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.301.03"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.117"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +16.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").Value = "2.666.57"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "2.308.40"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.807"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.57%  "
$ws.Range("D19").Value = "42.974.40"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.65%  "
$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.52%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.56%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0698"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "1.995.26"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").Value = "2.531.26"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  +2.83%  "
